# The workbook tracks daily price records for "Acelga" at the Macroferia
# Regional de Talca. A new record was inserted at row 227 (pushing all
# subsequent records down by one row), representing a new weekly entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 227; this shifts rows 227:305 down to 228:306
# and automatically extends the sheet dimension to A1:R306.
$ws.Rows(227).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A227").Value = 5
$ws.Range("B227").Value = "Macroferia Regional de Talca"
$ws.Range("C227").Value = "Maule"
$ws.Range("D227").Value = 44809
$ws.Range("E227").Value = 7
$ws.Range("F227").Value = 100112009
$ws.Range("G227").Value = "Acelga"
$ws.Range("H227").Value = "Sin especificar"
$ws.Range("I227").Value = "Primera"
$ws.Range("J227").Value = 500
$ws.Range("K227").Value = 2500
$ws.Range("L227").Value = 2500
$ws.Range("M227").Value = 2500
$ws.Range("N227").Value = "$/docena de atados (4 kilos)"
$ws.Range("O227").Value = "Región del Maule"
$ws.Range("P227").Value = 625
$ws.Range("Q227").Value = 4
$ws.Range("R227").Value = "Hortaliza"

# Match the date number format used by the rest of column D (D226 is an
# existing, correctly-formatted date cell directly above the new row).
$ws.Range("D227").NumberFormat = $ws.Range("D226").NumberFormat
